$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.379.06'
$ws.Range("E2").Value = '  -3.65%  '
$ws.Range("D3").Value = '2.994.23'
$ws.Range("E3").Value = '  -2.90%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '''537.79'
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").Value = '''134.60'
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '2.992.28'
$ws.Range("E8").Value = '  -2.73%  '
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '''0.148'
$ws.Range("E10").Value = '  -4.88%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '''6.13'
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("E12").Value = '  -2.54%  '
$ws.Range("D13").Value = '''0.0000222'
$ws.Range("E13").Value = '  -2.12%  '
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("D15").Value = '3.476.13'
$ws.Range("E15").Value = '  -3.19%  '
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").Value = '61.339.86'
$ws.Range("E17").Value = '  -3.78%  '
$ws.Range("D18").Value = '2.994.14'
$ws.Range("E18").Value = '  -3.10%  '
$ws.Range("E19").Value = '  -1.20%  '
$ws.Range("D20").Value = '''465.70'
$ws.Range("E20").Value = '  -4.58%  '
$ws.Range("E21").Value = '  -2.04%  '
$ws.Range("E22").Value = '  -3.76%  '
$ws.Range("D23").Value = '''6.93'
$ws.Range("E23").Value = '  -3.65%  '
$ws.Range("D24").Value = '''80.07'
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("D25").Value = '''11.97'
$ws.Range("E25").Value = '  -2.19%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").Value = '''2.68'
$ws.Range("E27").Value = '  -1.74%  '
$ws.Range("D28").Value = '''7.77'
$ws.Range("E28").Value = '  -5.99%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("E31").Value = '  +2.75%  '
$ws.Range("E32").Value = '  -2.63%  '
$ws.Range("D33").Value = '''5.50'
$ws.Range("E33").Value = '  +0.99%  '
$ws.Range("B34").Value = 'Stacks'
$ws.Range("C34").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D34").Value = '''2.28'
$ws.Range("E34").Value = '  -5.67%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = '''54.72'
$ws.Range("E35").Value = '  -4.51%  '
$ws.Range("E36").Value = '  -2.69%  '
$ws.Range("D37").Value = '''451.23'
$ws.Range("E37").Value = '  -8.66%  '
$ws.Range("D38").Value = '3.165.02'
$ws.Range("E38").Value = '  -3.02%  '
$ws.Range("D39").Value = '''0.0787'
$ws.Range("E39").Value = '  -1.62%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '''0.120'
$ws.Range("E40").Value = '  +2.03%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '''0.0385'
$ws.Range("E41").Value = '  -4.12%  '
$ws.Range("D42").Value = '''8.13'
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("D43").Value = '''2.45'
$ws.Range("E43").Value = '  -7.73%  '
$ws.Range("D44").Value = '''26.92'
$ws.Range("E44").Value = '  +8.42%  '
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("E46").Value = '  -5.13%  '
$ws.Range("E47").Value = '  -3.61%  '
$ws.Range("D48").Value = '''119.01'
$ws.Range("E48").Value = '  -1.97%  '
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("D50").Value = '0.0₃0494'
$ws.Range("E50").Value = '  -8.12%  '
$ws.Range("E51").Value = '  +5.75%  '
